# Update countries & provincias Spain
#
# Applies:
#  1) Refreshed case numbers for several existing country rows.
#  2) Four countries (Pakistan, Egipto, Republica de Chipre, Suazilandia)
#     moved to a new position (just above Japon, Finlandia, Letonia and
#     Benin respectively) together with refreshed case numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Plain data refreshes (no row movement)
# ---------------------------------------------------------------------

# Estados Unidos (row 4)
$ws.Range("B4").Value = 990021
$ws.Range("C4").Value = 2861
$ws.Range("D4").Value = 118869
$ws.Range("E4").Value = 815653
$ws.Range("F4").Value = 15143
$ws.Range("G4").Value = 86
$ws.Range("H4").Value = 55499

# Canada (row 15)
$ws.Range("B15").Value = 47319
$ws.Range("C15").Value = 424
$ws.Range("D15").Value = 17846
$ws.Range("E15").Value = 26856
$ws.Range("F15").Value = 557
$ws.Range("G15").Value = 57
$ws.Range("H15").Value = 2617

# Noruega (row 44)
$ws.Range("B44").Value = 7554
$ws.Range("C44").Value = 27
$ws.Range("D44").Value = 32
$ws.Range("E44").Value = 7317
$ws.Range("F44").Value = 52
$ws.Range("G44").Value = 4
$ws.Range("H44").Value = 205

# Moldavia (row 58)
$ws.Range("B58").Value = 3481
$ws.Range("C58").Value = 73
$ws.Range("D58").Value = 925
$ws.Range("E58").Value = 2454
$ws.Range("F58").Value = 212
$ws.Range("G58").Value = 6
$ws.Range("H58").Value = 102

# Bulgaria (row 84)
$ws.Range("B84").Value = 1363
$ws.Range("C84").Value = 63
$ws.Range("D84").Value = 206
$ws.Range("E84").Value = 1099
$ws.Range("F84").Value = 41
$ws.Range("G84").Value = 2
$ws.Range("H84").Value = 58

# Jordania (row 112)
$ws.Range("B112").Value = 449
$ws.Range("C112").Value = 2
$ws.Range("D112").Value = 342
$ws.Range("E112").Value = 100
$ws.Range("F112").Value = 5
$ws.Range("G112").Value = 0
$ws.Range("H112").Value = 7

# Isla de Man (row 124)
$ws.Range("B124").Value = 308
$ws.Range("C124").Value = 0
$ws.Range("D124").Value = 247
$ws.Range("E124").Value = 41
$ws.Range("F124").Value = 22
$ws.Range("G124").Value = 2
$ws.Range("H124").Value = 20

# ---------------------------------------------------------------------
# 2) Countries that moved position (cut from their old row, refreshed,
#    and pasted into a new row just above another country)
# ---------------------------------------------------------------------

# Pakistan: delete old row (after Chile), insert new row before Japon
$ws.Rows(32).Delete()
$ws.Rows(30).Insert()
$ws.Range("A30").Value = "Pakistan"
$ws.Range("B30").Value = 13909
$ws.Range("C30").Value = 581
$ws.Range("D30").Value = 3029
$ws.Range("E30").Value = 10588
$ws.Range("F30").Value = 111
$ws.Range("G30").Value = 11
$ws.Range("H30").Value = 292

# Egipto: delete old row (after Sudafrica), insert new row before Finlandia
$ws.Rows(54).Delete()
$ws.Rows(52).Insert()
$ws.Range("A52").Value = "Egipto"
$ws.Range("B52").Value = 4782
$ws.Range("C52").Value = 248
$ws.Range("D52").Value = 1236
$ws.Range("E52").Value = 3209
$ws.Range("F52").Value = 0
$ws.Range("G52").Value = 20
$ws.Range("H52").Value = 337

# Republica de Chipre: delete old row (after Letonia), insert new row before Letonia
$ws.Rows(93).Delete()
$ws.Rows(92).Insert()
$ws.Range("A92").Value = "Republica de Chipre"
$ws.Range("B92").Value = 822
$ws.Range("C92").Value = 5
$ws.Range("D92").Value = 148
$ws.Range("E92").Value = 660
$ws.Range("F92").Value = 15
$ws.Range("G92").Value = 0
$ws.Range("H92").Value = 14

# Suazilandia: delete old row (after Libia), insert new row before Benin
$ws.Rows(164).Delete()
$ws.Rows(162).Insert()
$ws.Range("A162").Value = "Suazilandia"
$ws.Range("B162").Value = 65
$ws.Range("C162").Value = 6
$ws.Range("D162").Value = 10
$ws.Range("E162").Value = 54
$ws.Range("F162").Value = 0
$ws.Range("G162").Value = 0
$ws.Range("H162").Value = 1
